# Regenerate save_data: replace Strike# values (column G, labeled "K") with
# recalculated strikeout counts for each outing row (rows 3-62).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K (strikeout) values, in row order starting at row 3 through row 62.
$newK = @(6,2,2,1,3,1,0,0,0,2,1,0,1,1,3,1,1,0,1,0,1,1,0,1,2,0,0,2,1,0,1,1,3,1,0,2,2,1,0,1,3,1,3,1,1,1,0,2,2,2,1,1,1,0,1,3,2,2,2,1)

$startRow = 3
for ($i = 0; $i -lt $newK.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $newK[$i]
}
